# "added new profile script"
# - delete the now-unused "Test Case Steps" sheet
# - append a new "Profile43" test case row to "Test Cases"
# - tidy up a couple of row heights / the description column width / selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Drop the empty "Test Case Steps" sheet entirely.
$steps = $wb.Worksheets.Item("Test Case Steps")
[void]$steps.Delete()

# New test case row (row 44): TCID, JIRA ID, Description, Runmode, Results.
# Pull formatting from the existing rows first (paste-formats), then overwrite
# the values so styles land the same way the prior rows are styled.
[void]$ws.Range("A43:E43").Copy()
[void]$ws.Range("A44:E44").PasteSpecial(-4122)
[void]$ws.Range("A42").Copy()
[void]$ws.Range("A44").PasteSpecial(-4122)
[void]$ws.Range("C43").Copy()
[void]$ws.Range("B44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A44").Value = "Profile43"
$ws.Range("B44").Value = "OPQA-2940"
$ws.Range("C44").Value = "Verify that user is able to update his first name and last name  fields with max length count "
$ws.Range("D44").Value = "Y"

# Rows 7 and 14 no longer need their explicit 30pt height.
$ws.Rows.Item(7).AutoFit()
$ws.Rows.Item(14).AutoFit()

# Description column grew a bit wider.
$ws.Columns.Item(3).ColumnWidth = 109.6

# Selection / scroll position as left by the author.
[void]$ws.Range("J31").Select()
